$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: tweak the wording of the "interactive talk" bullet - drop
# "in the chat " from the last paragraph, keeping everything else (and the
# run's formatting) intact.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shBullets = $s2.Shapes.Item(2)
$trBullets = $shBullets.TextFrame.TextRange
$lastParaIndex = $trBullets.Paragraphs().Count
$lastPara = $trBullets.Paragraphs($lastParaIndex, 1)
$oldLen = $lastPara.Length

$dash = [char]0x2013
$quote = [char]0x2019
$newLastParaText = "An interactive talk - I" + $quote + "ll ask you what you think will happen at various times " + $dash + " respond when this happens"

# Insert the replacement text in front of the existing run (this keeps the
# original run's rPr on the freshly inserted text), then drop the old text
# that is now trailing after it.
$lastPara.InsertBefore($newLastParaText) | Out-Null
$newLen = $newLastParaText.Length
$refreshedPara = $trBullets.Paragraphs($lastParaIndex, 1)
$oldTail = $refreshedPara.Characters($newLen + 1, $oldLen)
$oldTail.Delete()

# ---------------------------------------------------------------------------
# Slide 6: the cache-line-address diagram - resize/reposition the "set ="
# label + its three arrow connectors, and simplify the label text to a
# single expression/run.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# EMU -> point helper that compensates for the host's single-precision
# point storage + truncate-on-save EMU conversion, so the saved EMU values
# land exactly where we want them.
function EmuToPt([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

# "Rectangle 5" - the "set = ..." label box.
$shLabel = $s6.Shapes.Item(4)
$shLabel.Left   = EmuToPt 2685327
$shLabel.Top    = EmuToPt 206199
$shLabel.Width  = EmuToPt 3264891
$shLabel.Height = EmuToPt 327380

$trLabel = $shLabel.TextFrame.TextRange
$trLabel.Delete()
$trLabel.InsertAfter("set = (address & 0x1FC0)>>6") | Out-Null

# "Straight Arrow Connector 6"
$shConn1 = $s6.Shapes.Item(5)
$shConn1.Left   = EmuToPt 5950218
$shConn1.Top    = EmuToPt 369889
$shConn1.Width  = EmuToPt 965066
$shConn1.Height = EmuToPt 431390

# "Straight Arrow Connector 10"
$shConn2 = $s6.Shapes.Item(6)
$shConn2.Left   = EmuToPt 5950218
$shConn2.Top    = EmuToPt 369889
$shConn2.Width  = EmuToPt 965066
$shConn2.Height = EmuToPt 1977583

# "Straight Arrow Connector 12"
$shConn3 = $s6.Shapes.Item(7)
$shConn3.Left   = EmuToPt 5950218
$shConn3.Top    = EmuToPt 369889
$shConn3.Width  = EmuToPt 965066
$shConn3.Height = EmuToPt 3667487
